$d = $word.ActiveDocument

$d.Content.Find.Execute("21×20=420", $true, $false, $false, $false, $false, $true, 1, $false, "78×59=4602", 2) | Out-Null
$d.Content.Find.Execute("37×19=703", $true, $false, $false, $false, $false, $true, 1, $false, "57×12=684", 2) | Out-Null
$d.Content.Find.Execute("82×89=7298", $true, $false, $false, $false, $false, $true, 1, $false, "87×34=2958", 2) | Out-Null
$d.Content.Find.Execute("47×36=1692", $true, $false, $false, $false, $false, $true, 1, $false, "52×59=3068", 2) | Out-Null
$d.Content.Find.Execute("72×14=1008", $true, $false, $false, $false, $false, $true, 1, $false, "16×33=528", 2) | Out-Null
$d.Content.Find.Execute("86×54=4644", $true, $false, $false, $false, $false, $true, 1, $false, "54×49=2646", 2) | Out-Null
$d.Content.Find.Execute("14×92=1288", $true, $false, $false, $false, $false, $true, 1, $false, "96×54=5184", 2) | Out-Null
$d.Content.Find.Execute("68×79=5372", $true, $false, $false, $false, $false, $true, 1, $false, "28×57=1596", 2) | Out-Null
$d.Content.Find.Execute("22×70=1540", $true, $false, $false, $false, $false, $true, 1, $false, "41×94=3854", 2) | Out-Null
$d.Content.Find.Execute("62×72=4464", $true, $false, $false, $false, $false, $true, 1, $false, "83×43=3569", 2) | Out-Null
$d.Content.Find.Execute("26×75=1950", $true, $false, $false, $false, $false, $true, 1, $false, "25×43=1075", 2) | Out-Null
$d.Content.Find.Execute("37×88=3256", $true, $false, $false, $false, $false, $true, 1, $false, "17×95=1615", 2) | Out-Null
$d.Content.Find.Execute("57×15=855", $true, $false, $false, $false, $false, $true, 1, $false, "24×35=840", 2) | Out-Null
$d.Content.Find.Execute("51×76=3876", $true, $false, $false, $false, $false, $true, 1, $false, "96×98=9408", 2) | Out-Null
$d.Content.Find.Execute("12×22=264", $true, $false, $false, $false, $false, $true, 1, $false, "67×16=1072", 2) | Out-Null
$d.Content.Find.Execute("98×53=5194", $true, $false, $false, $false, $false, $true, 1, $false, "93×73=6789", 2) | Out-Null
$d.Content.Find.Execute("14×69=966", $true, $false, $false, $false, $false, $true, 1, $false, "99×75=7425", 2) | Out-Null
$d.Content.Find.Execute("55×49=2695", $true, $false, $false, $false, $false, $true, 1, $false, "83×26=2158", 2) | Out-Null
$d.Content.Find.Execute("74×40=2960", $true, $false, $false, $false, $false, $true, 1, $false, "39×63=2457", 2) | Out-Null
$d.Content.Find.Execute("40×68=2720", $true, $false, $false, $false, $false, $true, 1, $false, "51×12=612", 2) | Out-Null
$d.Content.Find.Execute("29×91=2639", $true, $false, $false, $false, $false, $true, 1, $false, "73×28=2044", 2) | Out-Null
$d.Content.Find.Execute("51×89=4539", $true, $false, $false, $false, $false, $true, 1, $false, "60×37=2220", 2) | Out-Null
$d.Content.Find.Execute("58×38=2204", $true, $false, $false, $false, $false, $true, 1, $false, "68×39=2652", 2) | Out-Null
$d.Content.Find.Execute("89×50=4450", $true, $false, $false, $false, $false, $true, 1, $false, "81×99=8019", 2) | Out-Null
$d.Content.Find.Execute("55×48=2640", $true, $false, $false, $false, $false, $true, 1, $false, "79×11=869", 2) | Out-Null
